$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 1064.5
$ws.Range("I41").Value = 1197.4
$ws.Range("K41").Value = 1197.4
$ws.Range("M41").Value = -757.4000000000001

# Row 64
$ws.Range("H64").Value = 4722.222
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496

# Row 67
$ws.Range("H67").Value = 4722.222
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716

# Row 80
$ws.Range("H80").Value = 4699.6
$ws.Range("I80").Value = 10500
$ws.Range("J80").Value = 832.6667
$ws.Range("K80").Value = 31500
$ws.Range("L80").Value = 2498.0001
$ws.Range("M80").Value = -30502
$ws.Range("N80").Value = -4494.0001

# Row 83
$ws.Range("H83").Value = 4699.6
$ws.Range("I83").Value = 10500
$ws.Range("J83").Value = 832.6667
$ws.Range("K83").Value = 94500
$ws.Range("L83").Value = 7494.0003
$ws.Range("M83").Value = -89508
$ws.Range("N83").Value = -17478.0003

# Row 107
$ws.Range("H107").Value = 421.61905
$ws.Range("J107").Value = 341.75
$ws.Range("L107").Value = 341.75
$ws.Range("N107").Value = -4181.75

# Row 137
$ws.Range("H137").Value = 7108.4165
$ws.Range("I137").Value = 6600.2383
$ws.Range("K137").Value = 19800.7149
$ws.Range("M137").Value = -17250.7149

# Row 141
$ws.Range("H141").Value = 9268.158
$ws.Range("I141").Value = 7859.1333
$ws.Range("K141").Value = 23577.3999
$ws.Range("M141").Value = -18397.3999

$ws = $wb.Worksheets.Item("ARM")
# Row 18
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10644

# Row 63
$ws.Range("H63").Value = 6862
$ws.Range("I63").Value = 2367.5
$ws.Range("K63").Value = 2367.5
$ws.Range("M63").Value = -1681.5

# Row 66
$ws.Range("H66").Value = 6862
$ws.Range("I66").Value = 2367.5
$ws.Range("K66").Value = 11837.5
$ws.Range("M66").Value = -8405.5

# Row 97
$ws.Range("H97").Value = 1132.6757
$ws.Range("I97").Value = 858.8077
$ws.Range("J97").Value = 1780
$ws.Range("K97").Value = 858.8077
$ws.Range("L97").Value = 1780
$ws.Range("M97").Value = -362.8077
$ws.Range("N97").Value = -2772

# Row 132
$ws.Range("H132").Value = 4094.7727
$ws.Range("I132").Value = 4094.7727
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12284.3181
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9754.3181
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2595.75
$ws.Range("I20").Value = 3059.7896
$ws.Range("K20").Value = 3059.7896
$ws.Range("M20").Value = -2812.7896

# Row 99
$ws.Range("H99").Value = 2340.205
$ws.Range("I99").Value = 1861.9656
$ws.Range("J99").Value = 3727.1
$ws.Range("K99").Value = 1861.9656
$ws.Range("L99").Value = 3727.1
$ws.Range("M99").Value = -363.9656
$ws.Range("N99").Value = -6723.1

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 866143.5
$ws.Range("I31").Value = 2794.3125
$ws.Range("K31").Value = 2794.3125
$ws.Range("M31").Value = -2499.3125

# Row 34
$ws.Range("H34").Value = 866143.5
$ws.Range("I34").Value = 2794.3125
$ws.Range("K34").Value = 2794.3125
$ws.Range("M34").Value = -2592.3125

# Row 58
$ws.Range("H58").Value = 4691.25
$ws.Range("I58").Value = 3410.6667
$ws.Range("K58").Value = 3410.6667
$ws.Range("M58").Value = -3207.6667

# Row 62
$ws.Range("H62").Value = 2006579.4
$ws.Range("I62").Value = 2006579.4
$ws.Range("K62").Value = 2006579.4
$ws.Range("M62").Value = -2005955.4

# Row 65
$ws.Range("H65").Value = 2006579.4
$ws.Range("I65").Value = 2006579.4
$ws.Range("K65").Value = 10032897
$ws.Range("M65").Value = -10029777

# Row 132
$ws.Range("H132").Value = 7898.6
$ws.Range("I132").Value = 4170.3335
$ws.Range("J132").Value = 13491
$ws.Range("K132").Value = 12511.0005
$ws.Range("L132").Value = 40473
$ws.Range("M132").Value = -9981.000499999998
$ws.Range("N132").Value = -45533

# Row 134
$ws.Range("H134").Value = 2960.95
$ws.Range("I134").Value = 2643.1052
$ws.Range("J134").Value = 9000
$ws.Range("K134").Value = 7929.3156
$ws.Range("L134").Value = 27000
$ws.Range("M134").Value = -5394.3156
$ws.Range("N134").Value = -32070

# Row 136
$ws.Range("H136").Value = 4691.25
$ws.Range("I136").Value = 3410.6667
$ws.Range("K136").Value = 10232.0001
$ws.Range("M136").Value = -7682.000100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 51947820
$ws.Range("I4").Value = 30014046
$ws.Range("K4").Value = 90042138
$ws.Range("M4").Value = -90042026

# Row 15
$ws.Range("H15").Value = 229.41176
$ws.Range("I15").Value = 200
$ws.Range("K15").Value = 600
$ws.Range("M15").Value = -460

# Row 98
$ws.Range("H98").Value = 1000
$ws.Range("I98").Value = 1000
$ws.Range("K98").Value = 3000
$ws.Range("M98").Value = -1502

# Row 114
$ws.Range("H114").Value = 2556.75
$ws.Range("I114").Value = 2076
$ws.Range("K114").Value = 6228
$ws.Range("M114").Value = -2974

# Row 132
$ws.Range("H132").Value = 2140
$ws.Range("I132").Value = 1480
$ws.Range("K132").Value = 13320
$ws.Range("M132").Value = -10790

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3222.1428
$ws.Range("I102").Value = 2590.6924
$ws.Range("J102").Value = 4248.25
$ws.Range("K102").Value = 2590.6924
$ws.Range("L102").Value = 4248.25
$ws.Range("M102").Value = -968.6923999999999
$ws.Range("N102").Value = -7492.25

# Row 126
$ws.Range("H126").Value = 4385.5713
$ws.Range("I126").Value = 3771.1428
$ws.Range("K126").Value = 11313.4284
$ws.Range("M126").Value = -8843.4284

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1491.1428
$ws.Range("I22").Value = 1541.5555
$ws.Range("J22").Value = 1400.4
$ws.Range("K22").Value = 1541.5555
$ws.Range("L22").Value = 1400.4
$ws.Range("M22").Value = -1246.5555
$ws.Range("N22").Value = -1990.4

# Row 27
$ws.Range("H27").Value = 1491.1428
$ws.Range("I27").Value = 1541.5555
$ws.Range("J27").Value = 1400.4
$ws.Range("K27").Value = 1541.5555
$ws.Range("L27").Value = 1400.4
$ws.Range("M27").Value = -1434.5555
$ws.Range("N27").Value = -1614.4

# Row 43
$ws.Range("H43").Value = 43372.5
$ws.Range("J43").Value = 40326.668
$ws.Range("L43").Value = 40326.668
$ws.Range("N43").Value = -40712.668

# Row 82
$ws.Range("H82").Value = 1788.5
$ws.Range("I82").Value = 1788.5
$ws.Range("K82").Value = 1788.5
$ws.Range("M82").Value = -1427.5

# Row 85
$ws.Range("H85").Value = 1788.5
$ws.Range("I85").Value = 1788.5
$ws.Range("K85").Value = 1788.5
$ws.Range("M85").Value = -540.5

# Row 122
$ws.Range("H122").Value = 6191.0713
$ws.Range("I122").Value = 6005.467
$ws.Range("K122").Value = 18016.401
$ws.Range("M122").Value = -15566.401

# Row 132
$ws.Range("H132").Value = 1103308.6
$ws.Range("I132").Value = 169548
$ws.Range("J132").Value = 2503949.5
$ws.Range("K132").Value = 508644
$ws.Range("L132").Value = 7511848.5
$ws.Range("M132").Value = -506114
$ws.Range("N132").Value = -7516908.5

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# Row 122
$ws.Range("H122").Value = 2830.077
$ws.Range("I122").Value = 2777.2778
$ws.Range("J122").Value = 2948.875
$ws.Range("K122").Value = 8331.8334
$ws.Range("L122").Value = 8846.625
$ws.Range("M122").Value = -5881.8334
$ws.Range("N122").Value = -13746.625

# Row 132
$ws.Range("H132").Value = 1116711.4
$ws.Range("I132").Value = 6779.4
$ws.Range("K132").Value = 20338.2
$ws.Range("M132").Value = -17808.2

# Row 136
$ws.Range("H136").Value = 6340
$ws.Range("I136").Value = 6765.8887
$ws.Range("J136").Value = 5792.4287
$ws.Range("K136").Value = 20297.6661
$ws.Range("L136").Value = 17377.2861
$ws.Range("M136").Value = -17747.6661
$ws.Range("N136").Value = -22477.2861
